$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores plain text (mixing thousand-separated values like
# "61.512.75" with plain decimals like "404.35"), so pre-format each cell we are
# about to rewrite as Text ("@") to stop Excel from reinterpreting plain-decimal
# looking values (e.g. "404.35", "133.28") as numbers.
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D11", "D13", "D15", "D16", "D17", "D19", "D21", "D22", "D23", "D26", "D27", "D28", "D29", "D30", "D31", "D33", "D35", "D37", "D38", "D41", "D42", "D48", "D49")
foreach ($addr in $priceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "61.512.75"
$ws.Range("E2").Value = "  -2.11%  "

$ws.Range("D3").Value = "3.377.07"
$ws.Range("E3").Value = "  -2.70%  "

$ws.Range("E4").Value = "  +0.01%  "

$ws.Range("D5").Value = "404.35"
$ws.Range("E5").Value = "  -2.42%  "

$ws.Range("D6").Value = "133.28"
$ws.Range("E6").Value = "  +7.45%  "

$ws.Range("D7").Value = "0.587"
$ws.Range("E7").Value = "  -0.83%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("E9").Value = "  -2.05%  "

$ws.Range("E10").Value = "  -7.17%  "

$ws.Range("D11").Value = "42.25"
$ws.Range("E11").Value = "  +2.44%  "

$ws.Range("E12").Value = "  -1.58%  "

$ws.Range("D13").Value = "3.905.40"
$ws.Range("E13").Value = "  -2.64%  "

$ws.Range("E14").Value = "  -2.21%  "

$ws.Range("D15").Value = "19.72"
$ws.Range("E15").Value = "  -1.22%  "

$ws.Range("D16").Value = "3.378.75"
$ws.Range("E16").Value = "  -2.63%  "

$ws.Range("D17").Value = "61.483.51"
$ws.Range("E17").Value = "  -1.95%  "

$ws.Range("E18").Value = "  -1.85%  "

$ws.Range("D19").Value = "10.92"
$ws.Range("E19").Value = "  +0.39%  "

$ws.Range("E20").Value = "  -7.36%  "

$ws.Range("D21").Value = "3.19"
$ws.Range("E21").Value = "  -3.99%  "

$ws.Range("D22").Value = "84.99"
$ws.Range("E22").Value = "  +3.65%  "

$ws.Range("D23").Value = "315.90"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("E24").Value = "  -1.27%  "

$ws.Range("E25").Value = "  -1.73%  "

$ws.Range("D26").Value = "4.79"
$ws.Range("E26").Value = "  +11.02%  "

$ws.Range("D27").Value = "29.41"
$ws.Range("E27").Value = "  -4.66%  "

$ws.Range("D28").Value = "8.22"
$ws.Range("E28").Value = "  +4.99%  "

$ws.Range("D29").Value = "7.65"
$ws.Range("E29").Value = "  -1.84%  "

$ws.Range("D30").Value = "2.68"
$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("D31").Value = "0.172"
$ws.Range("E31").Value = "  -1.89%  "

$ws.Range("E32").Value = "  -1.31%  "

$ws.Range("D33").Value = "11.35"
$ws.Range("E33").Value = "  -1.81%  "

$ws.Range("E34").Value = "  -0.26%  "

$ws.Range("D35").Value = "41.57"
$ws.Range("E35").Value = "  -1.30%  "

$ws.Range("E36").Value = "  -2.76%  "

$ws.Range("D37").Value = "51.67"
$ws.Range("E37").Value = "  -1.04%  "

$ws.Range("D38").Value = "0.998"

$ws.Range("E39").Value = "  -3.17%  "

$ws.Range("E40").Value = "  -3.27%  "

$ws.Range("D41").Value = "138.71"
$ws.Range("E41").Value = "  +1.98%  "

$ws.Range("D42").Value = "1.97"
$ws.Range("E42").Value = "  -1.14%  "

$ws.Range("E43").Value = "  -1.13%  "

$ws.Range("E44").Value = "  +3.29%  "

$ws.Range("E45").Value = "  +2.04%  "

$ws.Range("E46").Value = "  -1.44%  "

$ws.Range("E47").Value = "  -1.56%  "

$ws.Range("D48").Value = "21.47"
$ws.Range("E48").Value = "  -2.44%  "

$ws.Range("D49").Value = "2.120.83"
$ws.Range("E49").Value = "  -3.97%  "

$ws.Range("E50").Value = "  -7.31%  "

$ws.Range("E51").Value = "  +0.98%  "
